$d = $word.ActiveDocument

$pairs = @(
    @{old = "47×63=2961"; new = "31×38=1178"},
    @{old = "22×83=1826"; new = "90×74=6660"},
    @{old = "91×46=4186"; new = "36×18=648"},
    @{old = "72×92=6624"; new = "49×38=1862"},
    @{old = "32×32=1024"; new = "47×56=2632"},
    @{old = "38×11=418";  new = "44×68=2992"},
    @{old = "96×29=2784"; new = "15×25=375"},
    @{old = "16×33=528";  new = "42×37=1554"},
    @{old = "22×66=1452"; new = "55×37=2035"},
    @{old = "76×61=4636"; new = "25×68=1700"},
    @{old = "92×18=1656"; new = "29×48=1392"},
    @{old = "75×64=4800"; new = "98×13=1274"},
    @{old = "14×83=1162"; new = "60×94=5640"},
    @{old = "85×28=2380"; new = "60×18=1080"},
    @{old = "41×55=2255"; new = "61×17=1037"},
    @{old = "58×17=986";  new = "85×52=4420"},
    @{old = "72×65=4680"; new = "98×89=8722"},
    @{old = "94×34=3196"; new = "16×45=720"},
    @{old = "28×83=2324"; new = "39×50=1950"},
    @{old = "34×95=3230"; new = "90×83=7470"},
    @{old = "94×54=5076"; new = "40×59=2360"},
    @{old = "31×54=1674"; new = "68×32=2176"},
    @{old = "94×93=8742"; new = "85×67=5695"},
    @{old = "76×29=2204"; new = "47×74=3478"},
    @{old = "80×22=1760"; new = "38×67=2546"}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
